$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @(0.00107958050586058,0.000308451573103023,0.954503392967304,0.00694016039481801,0.00663170882171499,0.991980259099321,0.000925354719309068,0.0067859346082665,0.0183528685996299,0.0067859346082665,0.991980259099321,0.0067859346082665,0.00740283775447255,0,0.976711906230722,0,0,0.00308451573103023,0.999074645280691,0.000308451573103023,0.000154225786551511,0.000462677359654534,0.000771128932757557)
$row3 = @(0.000308451573103023,0.00740283775447255,0.00123380629241209,0.000154225786551511,0.992751388032079,0.00185070943861814,0,0.00169648365206663,0.0721776681061073,0.993059839605182,0,0.00462677359654534,0,0.00154225786551511,0.00154225786551511,0,0.998920419494139,0.970234423195558,0,0,0.0106415792720543,0.000154225786551511,0)
$row4 = @(0.998457742134485,0,0.0157310302282542,0.992751388032079,0.000154225786551511,0.0040098704503393,0.998920419494139,0,0.0280690931523751,0.000154225786551511,0.000925354719309068,0,0.9898210980876,0,0.0121838371375694,0.998920419494139,0,0,0.000771128932757557,0.999691548426897,0.000154225786551511,0.992288710672424,0.998920419494139)
$row5 = @(0,0.992134484885873,0.0222085132634176,0,0.000462677359654534,0.000925354719309068,0.000154225786551511,0.991054904380012,0.865515114127082,0,0.00663170882171499,0.986890808143122,0.000154225786551511,0.998457742134485,0.0067859346082665,0.000925354719309068,0.000925354719309068,0.0228254164096237,0,0,0.986582356570018,0.0067859346082665,0.000154225786551511)

for ($c = 0; $c -lt 23; $c++) {
    $ws.Cells.Item(2, $c + 2).Value = $row2[$c]
    $ws.Cells.Item(3, $c + 2).Value = $row3[$c]
    $ws.Cells.Item(4, $c + 2).Value = $row4[$c]
    $ws.Cells.Item(5, $c + 2).Value = $row5[$c]
}
